$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column widths (characters) chosen so that, after the host's internal
# points/pixel rounding, the persisted OOXML <col width="..."/> values land
# on the intended 13.2 / 9.6 / 2.4 / 10.8 / 12 layout used by the "Pays"
# calendar-style header grid.
$ws.Columns.Item(1).ColumnWidth  = 12.333333333333334   # -> 13.2
$ws.Columns.Item(2).ColumnWidth  = 8.833333333333334    # -> 9.6
$ws.Columns.Item(3).ColumnWidth  = 1.5                  # -> 2.4
$ws.Columns.Item(4).ColumnWidth  = 1.5                  # -> 2.4
$ws.Columns.Item(5).ColumnWidth  = 10.0                 # -> 10.8
$ws.Columns.Item(6).ColumnWidth  = 1.5                  # -> 2.4
$ws.Columns.Item(7).ColumnWidth  = 1.5                  # -> 2.4
$ws.Columns.Item(8).ColumnWidth  = 12.333333333333334   # -> 13.2
$ws.Columns.Item(9).ColumnWidth  = 1.5                  # -> 2.4
$ws.Columns.Item(10).ColumnWidth = 1.5                  # -> 2.4
$ws.Columns.Item(11).ColumnWidth = 11.166666666666666   # -> 12
$ws.Columns.Item(12).ColumnWidth = 1.5                  # -> 2.4
$ws.Columns.Item(13).ColumnWidth = 1.5                  # -> 2.4
$ws.Columns.Item(14).ColumnWidth = 8.833333333333334    # -> 9.6
$ws.Columns.Item(15).ColumnWidth = 1.5                  # -> 2.4
$ws.Columns.Item(16).ColumnWidth = 1.5                  # -> 2.4
$ws.Columns.Item(17).ColumnWidth = 11.166666666666666   # -> 12
$ws.Columns.Item(18).ColumnWidth = 1.5                  # -> 2.4
$ws.Columns.Item(19).ColumnWidth = 1.5                  # -> 2.4
$ws.Columns.Item(20).ColumnWidth = 8.833333333333334    # -> 9.6

# New "Date\Type" label in A2 — starts the notebook-style legend row under
# the existing day-of-week header, adding the 9th shared string.
$ws.Range("A2").Value = "Date\Type"
